$d = $word.ActiveDocument

$replacements = @(
    @{old = "13×32=416"; new = "12×64=768"},
    @{old = "17×66=1122"; new = "42×97=4074"},
    @{old = "62×87=5394"; new = "29×34=986"},
    @{old = "44×94=4136"; new = "26×48=1248"},
    @{old = "28×46=1288"; new = "71×48=3408"},
    @{old = "19×22=418"; new = "15×42=630"},
    @{old = "27×34=918"; new = "78×32=2496"},
    @{old = "52×38=1976"; new = "71×19=1349"},
    @{old = "16×56=896"; new = "12×69=828"},
    @{old = "94×44=4136"; new = "38×42=1596"},
    @{old = "63×23=1449"; new = "24×50=1200"},
    @{old = "62×73=4526"; new = "65×14=910"},
    @{old = "48×30=1440"; new = "12×96=1152"},
    @{old = "92×69=6348"; new = "80×94=7520"},
    @{old = "66×55=3630"; new = "65×99=6435"},
    @{old = "89×18=1602"; new = "48×40=1920"},
    @{old = "65×23=1495"; new = "16×65=1040"},
    @{old = "86×42=3612"; new = "66×92=6072"},
    @{old = "61×59=3599"; new = "85×68=5780"},
    @{old = "61×91=5551"; new = "26×29=754"},
    @{old = "71×23=1633"; new = "45×13=585"},
    @{old = "37×55=2035"; new = "49×57=2793"},
    @{old = "88×56=4928"; new = "14×24=336"},
    @{old = "25×36=900"; new = "99×90=8910"},
    @{old = "47×90=4230"; new = "69×85=5865"}
)

foreach ($r in $replacements) {
    $found = $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false,
                                      $true, 1, $false, $r.new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $($r.old)"
    }
}
